# Removed fetal deaths before PNC encounter.
# Zero out the weekly probability of fetal death (column C, rows 2-8 ==
# gestational weeks 0-6, i.e. before the first prenatal-care encounter
# around week 4) on the "potential_preg_untrt" (severity = 0 / low) sheet.
# All the dependent sheets/cells (moderate + high severity rows on this
# sheet, and the treated-group "potential_preg_trt" sheet) recompute off
# of these same input cells via formulas, so touching just C2:C8 here is
# enough to ripple the rest of the workbook.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("potential_preg_untrt")
$ws.Range("C2:C8").Value = 0

# Leave the freshly-edited sheet selected/active, matching the workbook's
# new "active tab" + selection state.
$ws.Activate()
$ws.Range("C9").Select() | Out-Null
